$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colIndices = @(2,3,5,6,7,8,9,10,12)

$data = @(
        @(1.613639150473716,0.2138367518421944,0.1773710331547171,2.030069189181162,0.962735258905326,0.9826694022769118,0.925217811486295,0.04640777681273178,0.5190011523732068),
        @(1.482210005852494,0.1929591715591243,0.1776150178559455,2.031209179659172,0.9644451287904587,0.9901085917916106,0.9370343875054559,0.04658510441611696,0.5063326831075017),
        @(1.401770107735388,0.1800566429994888,0.1778023963871185,2.033210858862077,0.9664006177526971,0.9953260563165998,0.9449859340126494,0.04671323602662625,0.4987214761434302),
        @(1.369056507789821,0.1747777585060533,0.1778882320456585,2.034353476638131,0.9674243269584792,0.997615325319785,0.9484008949082003,0.04677030440794461,0.4956621334130631),
        @(1.363628486434493,0.1738999362830782,0.1779030582831727,2.034562940461711,0.967607986974258,0.9980052999259783,0.9489784826212819,0.04678007410284479,0.4951566924101485),
        @(1.401328650313076,0.1799855351047199,0.1778035155828643,2.033224945480882,0.9664135066112749,0.9953562701518024,0.9450312827960126,0.04671398600146048,0.4986800452669087),
        @(1.568269598452105,0.2066555234208067,0.1774473789635351,2.030191832547573,0.9631363193425955,0.9850994048354664,0.9291474443005221,0.04646493109883565,0.51459847116503),
        @(1.897644961528897,0.2582933738855218,0.1770458553141605,2.03459381513062,0.9639371361734845,0.9701545355965777,0.90354086918353,0.04612877719786823,0.5471341624105293),
        @(2.140827595198857,0.2958340209530377,0.1769302421803474,2.044171136420744,0.9689910868529665,0.9623445442376521,0.888131104920248,0.04597399543991898,0.5718354991007857),
        @(2.251711628696455,0.3128274162958462,0.1769162716173707,2.049912997169926,0.9722728737370545,0.9594840725986273,0.8818652383911214,0.04592348053428452,0.583244637549484),
        @(2.293736860309878,0.3192503368308905,0.1769165070401044,2.052287004846278,0.9736579307437552,0.9585007761562139,0.8795999724187595,0.0459072032291381,0.5875896021791505),
        @(2.284684401988102,0.3178675849135288,0.1769162110225952,2.051766830339545,0.9733532884495162,0.9587080979473086,0.8800830500566903,0.04591058215798682,0.5866527468876654),
        @(2.255168359093148,0.3133560766579535,0.1769161804010828,2.050104303174535,0.9723839652535418,0.9594011720381985,0.8816767167548321,0.04592208429278344,0.5836016091814429),
        @(2.237093559027073,0.3105910698784555,0.1769168804334438,2.049111978091148,0.9718087902397485,0.9598387204093086,0.8826668951516581,0.04592950077266167,0.5817358928137963),
        @(2.133586155399939,0.2947217693662765,0.1769319299656456,2.043823807881822,0.9687964747897837,0.9625454437893382,0.8885556084231609,0.04597769624987791,0.5710933351914065),
        @(2.070153059200322,0.2849649251897972,0.1769510344441461,2.040934807577059,0.9672009271972968,0.9643834992625671,0.8923590690067726,0.04601235160575357,0.5646084658849873),
        @(2.03369248291898,0.2793451381574812,0.1769656598272302,2.039403471298215,0.9663756487050108,0.9655058450556879,0.8946167296470122,0.04603415769750541,0.5608947781926332),
        @(2.021351802770141,0.2774410174098989,0.1769712372467342,2.038907357521438,0.9661120694385943,0.9658970302440082,0.8953931445598258,0.04604186289373402,0.5596401828076267),
        @(2.076903098950652,0.286004376579541,0.1769486245288263,2.041228853165947,0.9673612017184325,0.9641810903911789,0.8919469348430766,0.04600846868068587,0.5652971123344344),
        @(2.263836974248363,0.3146815448261577,0.1769160396503242,2.050587204067412,0.9726648085808165,0.9591948853022814,0.8812056973676903,0.04591862851837902,0.5844971371293894),
        @(2.386217430471561,0.3333531395073805,0.1769269390952601,2.057867611770433,0.9769610118711256,0.9565185792592814,0.8748123864444963,0.04587653001545533,0.5971885867313063),
        @(2.320882065480248,0.3233942247366315,0.1769181851847126,2.053875227464061,0.9745917719376109,0.9578935638388657,0.8781671143268639,0.04589748129232518,0.5904018959955977),
        @(2.073851381094698,0.2855344731488003,0.1769497027057607,2.041095511541599,0.9672884550422651,0.9642723950564118,0.8921330395747304,0.04601021828695906,0.5649857299391812),
        @(1.808328960936592,0.2443942400878996,0.1771228621508669,2.032291465139394,0.9629405974258418,0.9736422070600526,0.9098722713168783,0.04620348317339662,0.5381918644929584)
    )

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $colIndices[$j]
        $ws.Cells.Item($excelRow, $col).Value2 = $row[$j]
    }
}
